# Consolidate the passenger ("SoCDTtiNTY-psgr") and freight ("SoCDTtiNTY-frgt")
# cargo-type sheets into a single "SoCDTtiNTY" sheet, distinguishing the two
# data blocks by prefixing the vehicle-type row labels with "passenger"/"freight".

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")
$wsFrgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

# Copy the freight data block (rows 2-7, cols A-H) underneath the passenger
# data block (which occupies rows 1-7) so it lands in rows 8-13.
$wsFrgt.Range("A2:H7").Copy()
$wsPsgr.Range("A8").PasteSpecial(-4104)

# Relabel column A: passenger rows (2-7) get a "passenger " prefix, and the
# newly appended freight rows (8-13) get a "freight " prefix.
$passengerLabels = @("passenger LDVs", "passenger HDVs", "passenger aircraft", "passenger rail", "passenger ships", "passenger motorbikes")
for ($i = 0; $i -lt $passengerLabels.Length; $i++) {
    $row = 2 + $i
    $wsPsgr.Range("A$row").Value = $passengerLabels[$i]
}

$freightLabels = @("freight LDVs", "freight HDVs", "freight aircraft", "freight rail", "freight ships", "freight motorbikes")
for ($i = 0; $i -lt $freightLabels.Length; $i++) {
    $row = 8 + $i
    $wsPsgr.Range("A$row").Value = $freightLabels[$i]
}

# The freight sheet's data now lives in the consolidated sheet, so drop it.
$wsFrgt.Delete()

# Rename the remaining (passenger) sheet to the consolidated sheet name.
$wsPsgr.Name = "SoCDTtiNTY"

# Tidy up the view state: reset selection to A1 and restore "About" as the
# active sheet (matching the original active/selected state).
$wsPsgr.Range("A1").Select()
$wb.Worksheets.Item("About").Activate()

Write-Host "Consolidated sheets:"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Name
}
